{"js": "const replacements = [\n  [\"2023-07-25 Tuesday\", \"2023-07-26 Wednesday\"],\n  [\"51+3=\", \"92-37=\"],\n  [\"15-13=\", \"79+4=\"],\n  [\"50-22=\", \"50-37=\"],\n  [\"73-58=\", \"63-60=\"],\n  [\"60+11=\", \"83-47=\"],\n  [\"10+62=\", \"87-83=\"],\n  [\"58-23=\", \"89-63=\"],\n  [\"40+20=\", \"66-1=\"],\n  [\"13+67=\", \"38-17=\"],\n  [\"74-15=\", \"68-29=\"],\n  [\"42+38=\", \"43+21=\"],\n  [\"72-15=\", \"99-34=\"],\n  [\"67-10=\", \"30+60=\"],\n  [\"10+11=\", \"3+23=\"],\n  [\"87-65=\", \"19+4=\"],\n  [\"26-22=\", \"13+20=\"],\n  [\"5+37=\", \"38-25=\"],\n  [\"19-11=\", \"56+33=\"],\n  [\"67+18=\", \"85-21=\"],\n  [\"10+4=\", \"87-20=\"],\n  [\"89-15=\", \"98-76=\"],\n  [\"29+55=\", \"82+17=\"],\n  [\"36+11=\", \"6+43=\"],\n  [\"55+38=\", \"23+45=\"],\n  [\"94-67=\", \"55-37=\"],\n  [\"1+56=\", \"96-55=\"],\n  [\"70-5=\", \"49+8=\"],\n  [\"26+60=\", \"54-26=\"],\n  [\"70+8=\", \"49+34=\"],\n  [\"89+0=\", \"57-42=\"],\n  [\"4+20=\", \"85-49=\"],\n  [\"86-59=\", \"20-4=\"],\n  [\"98-78=\", \"74+3=\"],\n  [\"86-64=\", \"2+95=\"],\n  [\"39+5=\", \"31+12=\"],\n  [\"55-43=\", \"51-6=\"],\n  [\"67+13=\", \"35-28=\"],\n  [\"74-34=\", \"53+1=\"],\n  [\"1+60=\", \"47+40=\"],\n  [\"77-9=\", \"62-32=\"],\n  [\"30-5=\", \"67+0=\"],\n  [\"50+42=\", \"2+63=\"],\n  [\"82-14=\", \"64-6=\"],\n  [\"66-20=\", \"71-41=\"],\n  [\"59+19=\", \"65-27=\"],\n  [\"69-51=\", \"25+53=\"],\n  [\"96-81=\", \"4+26=\"],\n  [\"24+0=\", \"36-34=\"],\n  [\"90-80=\", \"26+3=\"],\n  [\"72-43=\", \"13+16=\"],\n  [\"89-46=\", \"90-61=\"],\n  [\"22+22=\", \"48-3=\"],\n  [\"31+39=\", \"83-58=\"],\n  [\"45-41=\", \"67-24=\"],\n  [\"18+33=\", \"32-20=\"],\n  [\"43+22=\", \"96-73=\"],\n  [\"45-7=\", \"14+47=\"],\n  [\"83-81=\", \"70-31=\"],\n  [\"30+18=\", \"70-60=\"],\n  [\"10+50=\", \"89-63=\"],\n  [\"60-12=\", \"62-50=\"],\n  [\"64-61=\", \"69-28=\"],\n  [\"82-31=\", \"41-32=\"],\n  [\"47+31=\", \"57-31=\"],\n  [\"48+46=\", \"17+45=\"],\n  [\"67-55=\", \"80-52=\"],\n  [\"15+16=\", \"12-6=\"],\n  [\"59+6=\", \"98-40=\"],\n  [\"17-12=\", \"38+5=\"],\n  [\"9+28=\", \"7+64=\"],\n  [\"84-29=\", \"49-17=\"],\n  [\"19-7=\", \"73-65=\"],\n  [\"62+16=\", \"72-69=\"],\n  [\"25+51=\", \"1+76=\"],\n  [\"51-32=\", \"17+11=\"],\n  [\"44-1=\", \"21+19=\"],\n  [\"22+72=\", \"0+91=\"],\n  [\"22+27=\", \"70+10=\"],\n  [\"71+2=\", \"82-76=\"],\n  [\"81+2=\", \"11+46=\"],\n  [\"83-78=\", \"32+41=\"],\n  [\"83-33=\", \"10+52=\"],\n  [\"94-30=\", \"14+47=\"],\n  [\"58+39=\", \"89-22=\"],\n  [\"19+40=\", \"38+40=\"],\n  [\"34+55=\", \"79-69=\"],\n  [\"51-5=\", \"85-58=\"],\n  [\"15+7=\", \"99-77=\"],\n  [\"57-26=\", \"86-39=\"],\n  [\"72-40=\", \"1+33=\"],\n  [\"0+67=\", \"40-16=\"],\n  [\"87-39=\", \"29-20=\"],\n  [\"75-30=\", \"0+0=\"],\n  [\"20+23=\", \"24+48=\"],\n  [\"65+11=\", \"19+33=\"],\n  [\"3+15=\", \"57-3=\"],\n  [\"85-80=\", \"6+3=\"],\n  [\"83-82=\", \"71-51=\"],\n  [\"78-74=\", \"77+12=\"],\n  [\"91-32=\", \"99-72=\"],\n];\n\nconst body = context.document.body;\nlet applied = 0;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error('Text not found: ' + oldText);\n  }\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n    applied++;\n  }\n  await context.sync();\n}\nreturn 'applied:' + applied;\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2023-07-25 Tuesday\", \"2023-07-26 Wednesday\"),\n    @(\"51+3=\", \"92-37=\"),\n    @(\"15-13=\", \"79+4=\"),\n    @(\"50-22=\", \"50-37=\"),\n    @(\"73-58=\", \"63-60=\"),\n    @(\"60+11=\", \"83-47=\"),\n    @(\"10+62=\", \"87-83=\"),\n    @(\"58-23=\", \"89-63=\"),\n    @(\"40+20=\", \"66-1=\"),\n    @(\"13+67=\", \"38-17=\"),\n    @(\"74-15=\", \"68-29=\"),\n    @(\"42+38=\", \"43+21=\"),\n    @(\"72-15=\", \"99-34=\"),\n    @(\"67-10=\", \"30+60=\"),\n    @(\"10+11=\", \"3+23=\"),\n    @(\"87-65=\", \"19+4=\"),\n    @(\"26-22=\", \"13+20=\"),\n    @(\"5+37=\", \"38-25=\"),\n    @(\"19-11=\", \"56+33=\"),\n    @(\"67+18=\", \"85-21=\"),\n    @(\"10+4=\", \"87-20=\"),\n    @(\"89-15=\", \"98-76=\"),\n    @(\"29+55=\", \"82+17=\"),\n    @(\"36+11=\", \"6+43=\"),\n    @(\"55+38=\", \"23+45=\"),\n    @(\"94-67=\", \"55-37=\"),\n    @(\"1+56=\", \"96-55=\"),\n    @(\"70-5=\", \"49+8=\"),\n    @(\"26+60=\", \"54-26=\"),\n    @(\"70+8=\", \"49+34=\"),\n    @(\"89+0=\", \"57-42=\"),\n    @(\"4+20=\", \"85-49=\"),\n    @(\"86-59=\", \"20-4=\"),\n    @(\"98-78=\", \"74+3=\"),\n    @(\"86-64=\", \"2+95=\"),\n    @(\"39+5=\", \"31+12=\"),\n    @(\"55-43=\", \"51-6=\"),\n    @(\"67+13=\", \"35-28=\"),\n    @(\"74-34=\", \"53+1=\"),\n    @(\"1+60=\", \"47+40=\"),\n    @(\"77-9=\", \"62-32=\"),\n    @(\"30-5=\", \"67+0=\"),\n    @(\"50+42=\", \"2+63=\"),\n    @(\"82-14=\", \"64-6=\"),\n    @(\"66-20=\", \"71-41=\"),\n    @(\"59+19=\", \"65-27=\"),\n    @(\"69-51=\", \"25+53=\"),\n    @(\"96-81=\", \"4+26=\"),\n    @(\"24+0=\", \"36-34=\"),\n    @(\"90-80=\", \"26+3=\"),\n    @(\"72-43=\", \"13+16=\"),\n    @(\"89-46=\", \"90-61=\"),\n    @(\"22+22=\", \"48-3=\"),\n    @(\"31+39=\", \"83-58=\"),\n    @(\"45-41=\", \"67-24=\"),\n    @(\"18+33=\", \"32-20=\"),\n    @(\"43+22=\", \"96-73=\"),\n    @(\"45-7=\", \"14+47=\"),\n    @(\"83-81=\", \"70-31=\"),\n    @(\"30+18=\", \"70-60=\"),\n    @(\"10+50=\", \"89-63=\"),\n    @(\"60-12=\", \"62-50=\"),\n    @(\"64-61=\", \"69-28=\"),\n    @(\"82-31=\", \"41-32=\"),\n    @(\"47+31=\", \"57-31=\"),\n    @(\"48+46=\", \"17+45=\"),\n    @(\"67-55=\", \"80-52=\"),\n    @(\"15+16=\", \"12-6=\"),\n    @(\"59+6=\", \"98-40=\"),\n    @(\"17-12=\", \"38+5=\"),\n    @(\"9+28=\", \"7+64=\"),\n    @(\"84-29=\", \"49-17=\"),\n    @(\"19-7=\", \"73-65=\"),\n    @(\"62+16=\", \"72-69=\"),\n    @(\"25+51=\", \"1+76=\"),\n    @(\"51-32=\", \"17+11=\"),\n    @(\"44-1=\", \"21+19=\"),\n    @(\"22+72=\", \"0+91=\"),\n    @(\"22+27=\", \"70+10=\"),\n    @(\"71+2=\", \"82-76=\"),\n    @(\"81+2=\", \"11+46=\"),\n    @(\"83-78=\", \"32+41=\"),\n    @(\"83-33=\", \"10+52=\"),\n    @(\"94-30=\", \"14+47=\"),\n    @(\"58+39=\", \"89-22=\"),\n    @(\"19+40=\", \"38+40=\"),\n    @(\"34+55=\", \"79-69=\"),\n    @(\"51-5=\", \"85-58=\"),\n    @(\"15+7=\", \"99-77=\"),\n    @(\"57-26=\", \"86-39=\"),\n    @(\"72-40=\", \"1+33=\"),\n    @(\"0+67=\", \"40-16=\"),\n    @(\"87-39=\", \"29-20=\"),\n    @(\"75-30=\", \"0+0=\"),\n    @(\"20+23=\", \"24+48=\"),\n    @(\"65+11=\", \"19+33=\"),\n    @(\"3+15=\", \"57-3=\"),\n    @(\"85-80=\", \"6+3=\"),\n    @(\"83-82=\", \"71-51=\"),\n    @(\"78-74=\", \"77+12=\"),\n    @(\"91-32=\", \"99-72=\"),\n)\n\n$count = 0\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $ok = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $ok) {\n        throw \"Replacement failed for: $oldText\"\n    }\n    $count++\n}\n\nWrite-Output \"applied:$count\"\n"}
